$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.793.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.732.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '349.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.73%  '
$ws.Range("E7").Value = '  -2.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = '  -3.21%  '
$ws.Range("E10").Value = '  -3.29%  '
$ws.Range("E11").Value = '  +2.87%  '
$ws.Range("E12").Value = '  -2.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.32'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("E14").Value = '  -3.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.165.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.759.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.910'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '50.695.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.20%  '
$ws.Range("E20").Value = '  -3.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0946'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '261.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.51%  '
$ws.Range("E25").Value = '  -3.47%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.157'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("E30").Value = '  -1.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '51.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.95'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0436'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.91%  '
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.57%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.36%  '
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("E40").Value = '  -3.54%  '
$ws.Range("E41").Value = '  -2.87%  '
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.59%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.061.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("E48").Value = '  -1.92%  '
$ws.Range("B49").Value = 'SEI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.901'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.54%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.53%  '
